# Updates the cryptos price table (columns B-E) to the latest scraped
# values. Some Price (D) values look numeric ("1.004", "0.00000000350",
# "6.400" with a trailing zero, etc.) so a plain .Value assignment would
# let Excel silently coerce them into doubles (losing the trailing zero /
# exact text). Prefixing with a literal apostrophe forces text entry like
# typing it in the UI, and re-asserting Style "Normal" afterwards strips
# the quote-prefix formatting flag Excel tags the cell with so the cell
# style stays identical to the original (un-styled) cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.404.83'
$ws.Range("E2").Value = '  +0.64%  '

$ws.Range("D3").Value = '2.112.33'
$ws.Range("E3").Value = '  +2.31%  '

$ws.Range("E4").Value = '  +0.35%  '

$ws.Range("D5").Value = '''334.45'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.64%  '

$ws.Range("D6").Value = '''1.004'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").Value = '''0.5239'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.57%  '

$ws.Range("D8").Value = '''0.4555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.92%  '

$ws.Range("D9").Value = '''53.41'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +17.68%  '

$ws.Range("D10").Value = '''0.08918'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.47%  '

$ws.Range("D11").Value = '''1.178'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.76%  '

$ws.Range("D12").Value = '''24.42'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.89%  '

$ws.Range("D13").Value = '2.101.22'
$ws.Range("E13").Value = '  +2.02%  '

$ws.Range("D14").Value = '''6.845'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.63%  '

$ws.Range("D15").Value = '''8.020'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.54%  '

$ws.Range("D16").Value = '''96.70'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.40%  '

$ws.Range("D17").Value = '''1.005'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.16%  '

$ws.Range("E18").Value = '  +2.17%  '

$ws.Range("D19").Value = '''0.06648'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.71%  '

$ws.Range("D20").Value = '''19.28'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.49%  '

$ws.Range("D21").Value = '''1.003'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("D22").Value = '''6.369'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.97%  '

$ws.Range("D23").Value = '30.476.01'
$ws.Range("E23").Value = '  +0.80%  '

$ws.Range("D24").Value = '''12.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.77%  '

$ws.Range("E25").Value = '  +4.06%  '

$ws.Range("D26").Value = '2.351.00'
$ws.Range("E26").Value = '  +2.14%  '

$ws.Range("D27").Value = '''22.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '

$ws.Range("D28").Value = '''2.565'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.36%  '

$ws.Range("D29").Value = '''163.92'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.86%  '

$ws.Range("D30").Value = '''132.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.25%  '

$ws.Range("D31").Value = '''1.240'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.46%  '

$ws.Range("B32").Value = 'Stellar'
$ws.Range("C32").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D32").Value = '''0.1076'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.67%  '

$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").Value = '''1.704'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +14.37%  '

$ws.Range("D34").Value = '''6.400'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.07%  '

$ws.Range("D35").Value = '''3.933'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.64%  '

$ws.Range("D36").Value = '''10.52'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.50%  '

$ws.Range("D37").Value = '''0.02590'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.78%  '

$ws.Range("D38").Value = '''5.620'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.40%  '

$ws.Range("E39").Value = '  +4.36%  '

$ws.Range("D40").Value = '''0.2313'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.38%  '

$ws.Range("E41").Value = '  +2.72%  '

$ws.Range("D42").Value = '''0.6922'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.69%  '

$ws.Range("D43").Value = '''1.249'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.61%  '

$ws.Range("D44").Value = '''2.359'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.63%  '

$ws.Range("D45").Value = '''1.003'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.14%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''14.12'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.08%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = '''0.6402'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.50%  '

$ws.Range("D48").Value = '''3.662'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.97%  '

$ws.Range("D49").Value = '''0.00000000350'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +24.15%  '

$ws.Range("D50").Value = '''1.251'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.85%  '

$ws.Range("D51").Value = '''0.3418'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +27.17%  '
